$wb = $excel.ActiveWorkbook

# --- 1. Rename header cells on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add new "PO Forecast" sheet after "Monthly Trend" (last sheet) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the sheet-level layout options used on the other sheets
$wsForecast.Outline.SummaryRow = 1
$wsForecast.Outline.SummaryColumn = 1
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# --- 3. Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the header formatting used on the other sheets (bold, centered, bordered)
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# --- 4. Data rows ---
$wsForecast.Range("A2").Value = 44955.99999999999
$wsForecast.Range("B2").Value = 354
$wsForecast.Range("C2").Value = -208.2286980937472
$wsForecast.Range("D2").Value = 868.6312420980323
$wsForecast.Range("A3").Value = 44962.99999999999
$wsForecast.Range("B3").Value = 351
$wsForecast.Range("C3").Value = -187.1672404535002
$wsForecast.Range("D3").Value = 899.0171705489258
$wsForecast.Range("A4").Value = 44969.99999999999
$wsForecast.Range("B4").Value = 347
$wsForecast.Range("C4").Value = -178.0778418750713
$wsForecast.Range("D4").Value = 911.6868974689785
$wsForecast.Range("A5").Value = 44976.99999999999
$wsForecast.Range("B5").Value = 344
$wsForecast.Range("C5").Value = -195.1542428914131
$wsForecast.Range("D5").Value = 881.7218199648644
$wsForecast.Range("A6").Value = 45004.99999999999
$wsForecast.Range("B6").Value = 331
$wsForecast.Range("C6").Value = -184.8471530849855
$wsForecast.Range("D6").Value = 862.6603246326122
$wsForecast.Range("A7").Value = 45011.99999999999
$wsForecast.Range("B7").Value = 328
$wsForecast.Range("C7").Value = -184.1579277638263
$wsForecast.Range("D7").Value = 831.01626703895
$wsForecast.Range("A8").Value = 45018.99999999999
$wsForecast.Range("B8").Value = 325
$wsForecast.Range("C8").Value = -196.7439576609266
$wsForecast.Range("D8").Value = 875.1104932552884
$wsForecast.Range("A9").Value = 45025.99999999999
$wsForecast.Range("B9").Value = 322
$wsForecast.Range("C9").Value = -201.6306804453953
$wsForecast.Range("D9").Value = 840.806588575429
$wsForecast.Range("A10").Value = 45039.99999999999
$wsForecast.Range("B10").Value = 315
$wsForecast.Range("C10").Value = -208.0951435013897
$wsForecast.Range("D10").Value = 836.9955258056837
$wsForecast.Range("A11").Value = 45053.99999999999
$wsForecast.Range("B11").Value = 309
$wsForecast.Range("C11").Value = -285.2645484812236
$wsForecast.Range("D11").Value = 868.7328635878393
$wsForecast.Range("A12").Value = 45060.99999999999
$wsForecast.Range("B12").Value = 306
$wsForecast.Range("C12").Value = -233.7726793957667
$wsForecast.Range("D12").Value = 843.4119671615338
$wsForecast.Range("A13").Value = 45067.99999999999
$wsForecast.Range("B13").Value = 303
$wsForecast.Range("C13").Value = -242.2913035381042
$wsForecast.Range("D13").Value = 846.6758932226683
$wsForecast.Range("A14").Value = 45081.99999999999
$wsForecast.Range("B14").Value = 296
$wsForecast.Range("C14").Value = -200.5449423745298
$wsForecast.Range("D14").Value = 830.8418854985707
$wsForecast.Range("A15").Value = 45088.99999999999
$wsForecast.Range("B15").Value = 293
$wsForecast.Range("C15").Value = -229.0634337766708
$wsForecast.Range("D15").Value = 807.8461073371377
$wsForecast.Range("A16").Value = 45102.99999999999
$wsForecast.Range("B16").Value = 287
$wsForecast.Range("C16").Value = -235.2630168877782
$wsForecast.Range("D16").Value = 802.6684052008096
$wsForecast.Range("A17").Value = 45109.99999999999
$wsForecast.Range("B17").Value = 283
$wsForecast.Range("C17").Value = -218.2560847179437
$wsForecast.Range("D17").Value = 799.5891463018836
$wsForecast.Range("A18").Value = 45116.99999999999
$wsForecast.Range("B18").Value = 280
$wsForecast.Range("C18").Value = -232.3248514866688
$wsForecast.Range("D18").Value = 820.3218075318283
$wsForecast.Range("A19").Value = 45123.99999999999
$wsForecast.Range("B19").Value = 277
$wsForecast.Range("C19").Value = -246.141481440161
$wsForecast.Range("D19").Value = 822.6901643669639
$wsForecast.Range("A20").Value = 45130.99999999999
$wsForecast.Range("B20").Value = 274
$wsForecast.Range("C20").Value = -231.89499220334
$wsForecast.Range("D20").Value = 823.3090892444192
$wsForecast.Range("A21").Value = 45137.99999999999
$wsForecast.Range("B21").Value = 271
$wsForecast.Range("C21").Value = -276.1697007748388
$wsForecast.Range("D21").Value = 796.3946792781634
$wsForecast.Range("A22").Value = 45158.99999999999
$wsForecast.Range("B22").Value = 261
$wsForecast.Range("C22").Value = -270.3877319789033
$wsForecast.Range("D22").Value = 776.3743065850648
$wsForecast.Range("A23").Value = 45165.99999999999
$wsForecast.Range("B23").Value = 258
$wsForecast.Range("C23").Value = -244.3863891844487
$wsForecast.Range("D23").Value = 773.0884923979457
$wsForecast.Range("A24").Value = 45200.99999999999
$wsForecast.Range("B24").Value = 242
$wsForecast.Range("C24").Value = -265.5715675495114
$wsForecast.Range("D24").Value = 803.52099804616
$wsForecast.Range("A25").Value = 45221.99999999999
$wsForecast.Range("B25").Value = 232
$wsForecast.Range("C25").Value = -298.5628862439069
$wsForecast.Range("D25").Value = 747.5962385080012
$wsForecast.Range("A26").Value = 45228.99999999999
$wsForecast.Range("B26").Value = 229
$wsForecast.Range("C26").Value = -328.1748673058801
$wsForecast.Range("D26").Value = 742.6332876696666
$wsForecast.Range("A27").Value = 45235.99999999999
$wsForecast.Range("B27").Value = 226
$wsForecast.Range("C27").Value = -313.0076030689165
$wsForecast.Range("D27").Value = 744.6877463809692
$wsForecast.Range("A28").Value = 45242.99999999999
$wsForecast.Range("B28").Value = 222
$wsForecast.Range("C28").Value = -297.7564900081931
$wsForecast.Range("D28").Value = 774.5643968563959
$wsForecast.Range("A29").Value = 45249.99999999999
$wsForecast.Range("B29").Value = 219
$wsForecast.Range("C29").Value = -346.9706668144178
$wsForecast.Range("D29").Value = 745.8075560074816
$wsForecast.Range("A30").Value = 45256.99999999999
$wsForecast.Range("B30").Value = 216
$wsForecast.Range("C30").Value = -315.6527605634158
$wsForecast.Range("D30").Value = 712.6522609359782
$wsForecast.Range("A31").Value = 45263.99999999999
$wsForecast.Range("B31").Value = 213
$wsForecast.Range("C31").Value = -336.0619102992006
$wsForecast.Range("D31").Value = 692.5344235247335
$wsForecast.Range("A32").Value = 45270.99999999999
$wsForecast.Range("B32").Value = 210
$wsForecast.Range("C32").Value = -314.1942586523954
$wsForecast.Range("D32").Value = 741.8145929054431
$wsForecast.Range("A33").Value = 45277.99999999999
$wsForecast.Range("B33").Value = 206
$wsForecast.Range("C33").Value = -321.1811488833739
$wsForecast.Range("D33").Value = 744.577206806312
$wsForecast.Range("A34").Value = 45284.99999999999
$wsForecast.Range("B34").Value = 203
$wsForecast.Range("C34").Value = -312.5700264369655
$wsForecast.Range("D34").Value = 763.9408292293723
$wsForecast.Range("A35").Value = 45291.99999999999
$wsForecast.Range("B35").Value = 200
$wsForecast.Range("C35").Value = -314.6993022081697
$wsForecast.Range("D35").Value = 725.3057276706595
$wsForecast.Range("A36").Value = 45298.99999999999
$wsForecast.Range("B36").Value = 197
$wsForecast.Range("C36").Value = -338.3580070340769
$wsForecast.Range("D36").Value = 716.5363566912027
$wsForecast.Range("A37").Value = 45305.99999999999
$wsForecast.Range("B37").Value = 194
$wsForecast.Range("C37").Value = -418.5793511659029
$wsForecast.Range("D37").Value = 700.8416347313912
$wsForecast.Range("A38").Value = 45312.99999999999
$wsForecast.Range("B38").Value = 190
$wsForecast.Range("C38").Value = -349.0242475536018
$wsForecast.Range("D38").Value = 703.4539735593413
$wsForecast.Range("A39").Value = 45319.99999999999
$wsForecast.Range("B39").Value = 187
$wsForecast.Range("C39").Value = -362.3196754850358
$wsForecast.Range("D39").Value = 695.1373876945413

# Match the date-column formatting used on the other sheets
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A39").PasteSpecial(-4122)

$wsForecast.Range("A1").Select() | Out-Null
